$d = $word.ActiveDocument

# --- Paragraph 1: title block (date + line break + paper title) ---
$d.Content.Find.Execute("08.08.25", $true, $false, $false, $false, $false, $true, 1, $false, "06.08.25", 2) | Out-Null
$d.Content.Find.Execute("Efficient Attention Mechanisms for Large Language Models: A Survey", $true, $false, $false, $false, $false, $true, 1, $false, "Where to show Demos in Your Prompt: A Positional Bias of In-Context Learning", 2) | Out-Null

# --- Paragraph 2: intro paragraph ---
$d.Paragraphs.Item(2).Range.Text = "מאמר שנסקור היום מראה ששינוי פשוט במיקום הדוגמאות בפרומפט יכול לשנות דרמטית את רמת הדיוק של המודל. הנה מבט מהיר על הכלל החבוי הזה באינטראקציה עם בינה מלאכותית. מהנדסי פרומפטים אובססיביים לגבי התוכן של הפרומפטים שלהם. אבל המחברים חושפים שהתעלמנו ממשתנה קריטי לא פחות: המיקום של אותן דוגמאות. המחקר הזה לוקח את התחום מעבר למשחקי ניסוי וטעייה אל עבר מדע קפדני, והחידוש בו טמון בדיוק ובגישה השיטתית שלו."

# --- Paragraph 3: replaces old research-efforts paragraph ---
$d.Paragraphs.Item(3).Range.Text = "אמנם ידוע שהסדר הפנימי של דוגמאות משנה, אך מאמר זה מציג הבחנה מכרעת: לא מדובר בערבוב הדוגמאות, אלא בהזזת כל גוש הדוגמאות, ללא שינוי, למיקומים מבניים שונים בתוך הפרומפט. המחברים מכנים תופעה ספציפית זו הטיית DPP (DEMOS POSITION IN PROMPT). כדי לחקור זאת, הם יצרו מסגרת שיטתית הבוחנת ארבעה מיקומים קנוניים: בתחילת או בסוף הנחיות המערכת, ובתחילת או בסוף הודעת המשתמש. גישה זו הופכת תצפית מעורפלת למדע שניתן לבחון."

# --- Delete the old "Heading4" paragraph (4 משפחות היעילות) ---
$d.Paragraphs.Item(4).Range.Delete()

# --- After deletion, old paragraphs 5..10 shifted down to 4..9 ---
$d.Paragraphs.Item(4).Range.Text = "המחברים מסתכלים מעבר לדיוק פשוט על ידי מדידת PREDICTION-CHANGE המודד כמה תשובות בפועל מתהפכות כאשר מבנה הפרומפט משתנה. זוהי תרומה חיונית, מכיוון שהיא חושפת חוסר יציבות סמוי. מודל עשוי להיראות מדויק באותה מידה עם שני פרומפטים שונים, אך אחד מהם עלול לגרום להתנהגות בלתי צפויה לחלוטין."
$d.Paragraphs.Item(5).Range.Text = "המחקר רחב ההיקף, שכלל עשרה מודלים ושמונה משימות שונות, הניב תוצאות ברורות וניתנות ליישום."
$d.Paragraphs.Item(6).Range.Text = "- אפקט הראשוניות הוא אמיתי: מיקום דוגמאות בשלב מוקדם בפרומפט (ssp, esp) מניב באופן עקבי דיוק גבוה יותר ויציבות רבה יותר, עם שיפור של עד 6 נקודות דיוק."

# --- Paragraph (old #8, now #7) carried an xml:space="preserve" <w:t>; use Find so the
#     stale attribute is not inherited by Range.Text assignment ---
$d.Content.Find.Execute("משפחה זו כוללת שיטות המשתמשות sliding windows, שבהן טוקן מתייחס רק לשכניו המקומיים. גישה זו מבוססת על האינטואיציה החזקה של ""מקומיות ההקשר"" (locality of reference) שמילים סמוכות הן לרוב הרלוונטיות ביותר. כדי למנוע אובדן של מידע גלובלי, גישה זו מחוזקת לעיתים קרובות באמצעות מספר טוקנים גלובליים הרשאים להתייחס לכל הרצף, או באמצעות תבניות מורחבות/מדלגות (dilated/strided patterns) המדלגות באופן שיטתי על טוקנים כדי לכסות שדה קליטה רחב יותר עם מספר קבוע של חישובים. ", $true, $false, $false, $false, $false, $true, 1, $false, "- אזור הסכנה: הצבת דוגמאות בסוף (eum) היא לרוב הרסנית. היא גורמת לירידה משמעותית בביצועים ולתנודתיות גבוהה, והופכת מעל 30% מהתשובות של המודל במשימות מסוימות של שאלות ותשובות, מבלי לשפר את נכונותן.", 2) | Out-Null

$d.Paragraphs.Item(8).Range.Text = "- אין פתרון קסם: המיקום האופטימלי אינו אוניברסלי; הוא תלוי בגודל המודל ובסוג המשימה. לדוגמה, בעוד שמודלים קטנים יותר מעדיפים דוגמאות בתחילה, מודל גדול כמו LLAMA3-70B מעדיף לעיתים קרובות שהדוגמאות יהיו קרובות יותר לשאילתה (sum)."
$d.Paragraphs.Item(9).Range.Text = "המחקר מבהיר: מיקום הדוגמאות שלכם אינו בחירה סגנונית. זהו פרמטר קריטי שיש לבחון ולהתאים. הסתמכות על פורמט ברירת מחדל עלולה לבזבז ביצועים ויציבות משמעותיים. לראשונה, ישנה מפת דרכים ברורה להבנה ואופטימיזציה של המימד החיוני הזה בעיצוב פרומפטים."

# --- Delete old paragraphs 10..18 (the remainder of the 4-families deep dive) ---
$pStart = $d.Paragraphs.Item(10)
$pEnd = $d.Paragraphs.Item(18)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.Delete()

# --- Final paragraph: arxiv link ---
$d.Paragraphs.Item(10).Range.Text = "https://arxiv.org/abs/2507.22887"

